# Applies the edits described by the commit diff to Project_Observations.docx.
$d = $word.ActiveDocument

# (2) Item 2: "... project state." -> "... project state with the provided data."
$d.Content.Find.Execute(
    "the ratio of successful projects is more than any other project state.",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "the ratio of successful projects is more than any other project state with the provided data.",
    2)

# (3 setup happens later, after text reflows) -- bookmark move is done last below.

# (4) "Limitations of the data set" paragraph rewritten
$d.Content.Find.Execute(
    "The data set is huge and calculations using excel may result in unexpected results and its hard to figure out these issues or errors.",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "The data set provided is not a complete data set and so analysis might not be accurate. Also, calculations using excel may result in unexpected results and it is hard to figure out these issues or errors.",
    2)

# (5) Typo fix: "int his" -> "in this"
$d.Content.Find.Execute(
    "excel int his worksheet",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "excel in this worksheet",
    2)

# (6) "These look pretty organized and structured when" -> "These look organized and structured when"
$d.Content.Find.Execute(
    "These look pretty organized and structured when",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "These look organized and structured when",
    2)

# (7) ", column (state) and field value(count of state)" -> ", column (state) and field value (count of state)"
$d.Content.Find.Execute(
    "column (state) and field value(count of state)",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "column (state) and field value (count of state)",
    2)

# (3) Move the "_GoBack" bookmark from right after item 4) ("...outcomes sheet.")
#     to sit just before "(count of state)" in the pivot-table bullet item.
#     Adding a bookmark with the same name ("_GoBack") re-seats/moves the
#     existing bookmark rather than creating a duplicate.
$bmRange = $d.Content
$bmRange.Find.Execute(
    "field value ",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "", 0)
$bmRange.Collapse(0)
$d.Bookmarks.Add("_GoBack", $bmRange)
